# Fruta / hortaliza, semanal
# Row 2 and Row 3 data (Fecha + price columns) are swapped between each other,
# matching the weekly price-report logic: the row that used to report the
# 44200 reading now reports the 44210 reading, and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 becomes the old Row 3 values ---
$ws.Range("D2").Value = 44210
$ws.Range("J2").Value = 1450
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = 1650
$ws.Range("P2").Value = 1650

# --- Row 3 becomes the old Row 2 values ---
$ws.Range("D3").Value = 44200
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1450
$ws.Range("P3").Value = 1450
